$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed invoice data (rows 2-7). Row 1 is the header: #, ID, Due Date, Invoice Url.
# Columns A (#) and C (Due Date) are numeric/date-looking text, so each such
# cell is forced to Text format before the value is written - otherwise Excel
# would silently reinterpret "10" as a number or "06-08-2024" as a date,
# which would break downstream string comparisons in the RPA flow.
$data = @(
    @{ Row = 2; Num = "1";  Id = "gw2odrctzelzpzjf3mmslp"; Due = "28-07-2024"; Url = "https://rpachallengeocr.azurewebsites.net/invoices/12.jpg" },
    @{ Row = 3; Num = "3";  Id = "538o6gypasvtnt0sfq1ivr"; Due = "06-08-2024"; Url = "https://rpachallengeocr.azurewebsites.net/invoices/8.jpg" },
    @{ Row = 4; Num = "9";  Id = "zs2l4nrxtyprw8hh31dpdq"; Due = "24-07-2024"; Url = "https://rpachallengeocr.azurewebsites.net/invoices/7.jpg" },
    @{ Row = 5; Num = "10"; Id = "ste1r42xtv20yf084aj2t";  Due = "20-08-2024"; Url = "https://rpachallengeocr.azurewebsites.net/invoices/11.jpg" },
    @{ Row = 6; Num = "11"; Id = "snspz98nxw80ux68n4q69n"; Due = "10-07-2024"; Url = "https://rpachallengeocr.azurewebsites.net/invoices/10.jpg" },
    @{ Row = 7; Num = "12"; Id = "y0yvyji9lgkln1b210be98"; Due = "19-07-2024"; Url = "https://rpachallengeocr.azurewebsites.net/invoices/7.jpg" }
)

foreach ($row in $data) {
    $r = $row.Row

    $cellA = $ws.Cells.Item($r, 1)
    if ($cellA.Value2 -ne $row.Num) {
        $cellA.NumberFormat = "@"
        $cellA.Value = $row.Num
    }

    $cellB = $ws.Cells.Item($r, 2)
    if ($cellB.Value2 -ne $row.Id) {
        $cellB.Value = $row.Id
    }

    $cellC = $ws.Cells.Item($r, 3)
    if ($cellC.Value2 -ne $row.Due) {
        $cellC.NumberFormat = "@"
        $cellC.Value = $row.Due
    }

    $cellD = $ws.Cells.Item($r, 4)
    if ($cellD.Value2 -ne $row.Url) {
        $cellD.Value = $row.Url
    }
}
